$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 399
$arr[0,1] = 499
$arr[0,2] = 349
$arr[0,3] = 499
$arr[0,4] = 349
$arr[0,5] = -329
$arr[0,6] = -689
$ws.Range("H12:N12").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1443.6875
$arr[0,1] = 0
$arr[0,2] = 1443.6875
$arr[0,3] = 0
$arr[0,4] = 4331.0625
$arr[0,5] = $null
$arr[0,6] = -4667.0625
$ws.Range("H17:N17").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2560.4
$arr[0,1] = 1800
$arr[0,2] = 2750.5
$arr[0,3] = 1800
$arr[0,4] = 2750.5
$arr[0,5] = -1625
$arr[0,6] = -3100.5
$ws.Range("H40:N40").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3389.4375
$arr[0,1] = 2769.3333
$arr[0,2] = 5249.75
$arr[0,3] = 2769.3333
$arr[0,4] = 5249.75
$arr[0,5] = -2329.3333
$arr[0,6] = -6129.75
$ws.Range("H41:N41").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1375
$arr[0,1] = 1000
$arr[0,2] = 2500
$arr[0,3] = 1000
$arr[0,4] = 2500
$arr[0,5] = -931
$arr[0,6] = -2638
$ws.Range("H43:N43").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1478
$arr[0,1] = 1220
$arr[0,2] = 1650
$arr[0,3] = 1220
$arr[0,4] = 1650
$arr[0,5] = -596
$arr[0,6] = -2898
$ws.Range("H62:N62").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1478
$arr[0,1] = 1220
$arr[0,2] = 1650
$arr[0,3] = 6100
$arr[0,4] = 8250
$arr[0,5] = -2980
$arr[0,6] = -14490
$ws.Range("H65:N65").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3998.3333
$arr[0,1] = 0
$arr[0,2] = 3998.3333
$arr[0,3] = 0
$arr[0,4] = 3998.3333
$arr[0,5] = $null
$arr[0,6] = -4810.3333
$ws.Range("H88:N88").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3998.3333
$arr[0,1] = 0
$arr[0,2] = 3998.3333
$arr[0,3] = 0
$arr[0,4] = 3998.3333
$arr[0,5] = $null
$arr[0,6] = -6806.3333
$ws.Range("H91:N91").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1221
$arr[0,1] = 957
$arr[0,2] = 1749
$arr[0,3] = 2871
$arr[0,4] = 5247
$arr[0,5] = 196
$arr[0,6] = -11381
$ws.Range("H111:N111").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2111.111
$arr[0,1] = 2000
$arr[0,2] = 3000
$arr[0,3] = 2000
$arr[0,4] = 3000
$arr[0,5] = 1254
$arr[0,6] = -9508
$ws.Range("H113:N113").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H125:N125").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2738.5715
$arr[0,1] = 2738.5715
$arr[0,2] = 0
$arr[0,3] = 8215.7145
$arr[0,4] = 0
$arr[0,5] = -3175.7145
$arr[0,6] = $null
$ws.Range("H131:N131").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2740.1538
$arr[0,1] = 2768.1428
$arr[0,2] = 2707.5
$arr[0,3] = 8304.4284
$arr[0,4] = 8122.5
$arr[0,5] = -5754.428400000001
$arr[0,6] = -13222.5
$ws.Range("H137:N137").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2564.6667
$arr[0,1] = 2347
$arr[0,2] = 3000
$arr[0,3] = 7041
$arr[0,4] = 9000
$arr[0,5] = -1901
$arr[0,6] = -19280
$ws.Range("H138:N138").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 900.8
$arr[0,1] = 900.8
$arr[0,2] = 0
$arr[0,3] = 2702.4
$arr[0,4] = 0
$arr[0,5] = 2477.6
$arr[0,6] = $null
$ws.Range("H141:N141").Value2 = $arr

$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2250
$arr[0,1] = 2250
$arr[0,2] = 0
$arr[0,3] = 2250
$arr[0,4] = 0
$arr[0,5] = -1564
$arr[0,6] = $null
$ws.Range("H63:N63").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2250
$arr[0,1] = 2250
$arr[0,2] = 0
$arr[0,3] = 11250
$arr[0,4] = 0
$arr[0,5] = -7818
$arr[0,6] = $null
$ws.Range("H66:N66").Value2 = $arr

$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 33741.5
$arr[0,1] = 33741.5
$arr[0,2] = 0
$arr[0,3] = 33741.5
$arr[0,4] = 0
$arr[0,5] = -33257.5
$arr[0,6] = $null
$ws.Range("H54:N54").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1401.8667
$arr[0,1] = 1437.7142
$arr[0,2] = 900
$arr[0,3] = 1437.7142
$arr[0,4] = 900
$arr[0,5] = -314.7141999999999
$arr[0,6] = -3146
$ws.Range("H86:N86").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1401.8667
$arr[0,1] = 1437.7142
$arr[0,2] = 900
$arr[0,3] = 7188.571
$arr[0,4] = 4500
$arr[0,5] = -1572.571
$arr[0,6] = -15732
$ws.Range("H89:N89").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 18317
$arr[0,1] = 18317
$arr[0,2] = 0
$arr[0,3] = 18317
$arr[0,4] = 0
$arr[0,5] = -16570
$arr[0,6] = $null
$ws.Range("H105:N105").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2009.2
$arr[0,1] = 753.5
$arr[0,2] = 2846.3333
$arr[0,3] = 2260.5
$arr[0,4] = 8538.999899999999
$arr[0,5] = 274.5
$arr[0,6] = -13608.9999
$ws.Range("H134:N134").Value2 = $arr

$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 952.6667
$arr[0,1] = 1003.6667
$arr[0,2] = 927.1667
$arr[0,3] = 1003.6667
$arr[0,4] = 927.1667
$arr[0,5] = -716.6667
$arr[0,6] = -1501.1667
$ws.Range("H16:N16").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H62:N62").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H65:N65").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1914.1666
$arr[0,1] = 1997
$arr[0,2] = 1500
$arr[0,3] = 1997
$arr[0,4] = 1500
$arr[0,5] = -77
$arr[0,6] = -5340
$ws.Range("H107:N107").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 952.6667
$arr[0,1] = 1003.6667
$arr[0,2] = 927.1667
$arr[0,3] = 1003.6667
$arr[0,4] = 927.1667
$arr[0,5] = 1166.3333
$arr[0,6] = -5267.1667
$ws.Range("H113:N113").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5956.857
$arr[0,1] = 4424.75
$arr[0,2] = 7999.6665
$arr[0,3] = 13274.25
$arr[0,4] = 23998.9995
$arr[0,5] = -10744.25
$arr[0,6] = -29058.9995
$ws.Range("H132:N132").Value2 = $arr

$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 11.222222
$arr[0,1] = 13
$arr[0,2] = 6.6
$arr[0,3] = 78
$arr[0,4] = 39.59999999999999
$arr[0,5] = 35
$arr[0,6] = -265.6
$ws.Range("H2:N2").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3
$arr[0,1] = 0
$arr[0,2] = 3
$arr[0,3] = 0
$arr[0,4] = 9
$arr[0,5] = $null
$arr[0,6] = -233
$ws.Range("H7:N7").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 49
$arr[0,1] = 1
$arr[0,2] = 61
$arr[0,3] = 3
$arr[0,4] = 183
$arr[0,5] = 170
$arr[0,6] = -529
$ws.Range("H12:N12").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 329.33334
$arr[0,1] = 45.5
$arr[0,2] = 471.25
$arr[0,3] = 136.5
$arr[0,4] = 1413.75
$arr[0,5] = 98.5
$arr[0,6] = -1883.75
$ws.Range("H23:N23").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 123.57143
$arr[0,1] = 87.5
$arr[0,2] = 138
$arr[0,3] = 525
$arr[0,4] = 828
$arr[0,5] = -242
$arr[0,6] = -1394
$ws.Range("H33:N33").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 0
$arr[0,1] = 0
$arr[0,2] = 0
$arr[0,3] = 0
$arr[0,4] = 0
$arr[0,5] = $null
$arr[0,6] = $null
$ws.Range("H57:N57").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 650
$arr[0,1] = 650
$arr[0,2] = 0
$arr[0,3] = 1950
$arr[0,4] = 0
$arr[0,5] = 2140
$arr[0,6] = $null
$ws.Range("H110:N110").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 9000
$arr[0,1] = 9000
$arr[0,2] = 0
$arr[0,3] = 27000
$arr[0,4] = 0
$arr[0,5] = -21980
$arr[0,6] = $null
$ws.Range("H130:N130").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2745.375
$arr[0,1] = 4965
$arr[0,2] = 2005.5
$arr[0,3] = 14895
$arr[0,4] = 6016.5
$arr[0,5] = -9855
$arr[0,6] = -16096.5
$ws.Range("H131:N131").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1500
$arr[0,1] = 1000
$arr[0,2] = 2000
$arr[0,3] = 3000
$arr[0,4] = 6000
$arr[0,5] = 2140
$arr[0,6] = -16280
$ws.Range("H138:N138").Value2 = $arr

$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1500
$arr[0,1] = 0
$arr[0,2] = 1500
$arr[0,3] = 0
$arr[0,4] = 1500
$arr[0,5] = $null
$arr[0,6] = -2558
$ws.Range("H22:N22").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 5299.8335
$arr[0,1] = 5159.8
$arr[0,2] = 6000
$arr[0,3] = 5159.8
$arr[0,4] = 6000
$arr[0,5] = -5008.8
$arr[0,6] = -6302
$ws.Range("H43:N43").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1215.5
$arr[0,1] = 555
$arr[0,2] = 1435.6666
$arr[0,3] = 1665
$arr[0,4] = 4306.9998
$arr[0,5] = 785
$arr[0,6] = -9206.9998
$ws.Range("H122:N122").Value2 = $arr

$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 40500
$arr[0,1] = 0
$arr[0,2] = 40500
$arr[0,3] = 0
$arr[0,4] = 40500
$arr[0,5] = $null
$arr[0,6] = -40950
$ws.Range("H64:N64").Value2 = $arr

$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 40500
$arr[0,1] = 0
$arr[0,2] = 40500
$arr[0,3] = 0
$arr[0,4] = 40500
$arr[0,5] = $null
$arr[0,6] = -42060
$ws.Range("H67:N67").Value2 = $arr

$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 937.625
$arr[0,1] = 1108.5
$arr[0,2] = 425
$arr[0,3] = 2217
$arr[0,4] = 1000
$arr[0,5] = -1676
$arr[0,6] = -1932
$ws.Range("H100:N100").Value2 = $arr
